$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (bold, bordered style matching the rest of row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-32
$data = @{
    "2"  = @(5,6)
    "3"  = @(8,8)
    "4"  = @(6,6)
    "5"  = @(8,8)
    "6"  = @(7,7)
    "7"  = @(8,8)
    "8"  = @(7,7)
    "9"  = @(8,8)
    "10" = @(9,9)
    "11" = @(9,9)
    "12" = @(6,6)
    "13" = @(7,7)
    "14" = @(9,9)
    "15" = @(7,7)
    "16" = @(7,7)
    "17" = @(5,6)
    "18" = @(8,8)
    "19" = @(1,5)
    "20" = @(1,4)
    "21" = @(1,4)
    "22" = @(6,6)
    "23" = @(9,9)
    "24" = @(7,8)
    "25" = @(9,9)
    "26" = @(8,8)
    "27" = @(7,7)
    "28" = @(9,9)
    "29" = @(6,6)
    "30" = @(5,5)
    "31" = @(5,6)
    "32" = @(3,4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $r = [int]$row
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
